$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 2) to upper case column names
$ws.Range("A2").Value = "STATION"
$ws.Range("B2").Value = "NAME"
$ws.Range("C2").Value = "NETID"
$ws.Range("D2").Value = "LAT"
$ws.Range("E2").Value = "LON"
$ws.Range("F2").Value = "INTENSITY"

# Move the active selection to F2
$ws.Range("F2").Select()
